# Box Plot Updates, Color Updates Main Figures
#
# Repositions several shapes inside the chart group on slide 1:
#   - "pl8"  (a freeform line/leader): moves/resizes slightly
#   - "tx9".."tx16" (text labels): reposition (some only by a few EMU)
#
# NOTE: PowerPoint COM exposes shape position/size (Left/Top/Width/Height)
# as single-precision (float32) "points" values. 1 pt = 12700 EMU, and the
# round-trip float64(pt) -> float32 -> *12700 -> truncate-to-EMU can lose
# the final EMU of precision for some values. The literal point values
# below were chosen so that, after that float32 round-trip, they land on
# the exact target EMU offsets from the authoritative OOXML.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$g = $s.Shapes.Item(2)

# "pl8" - small freeform connector/leader line: off x=5225794 y=2766416 (EMU),
# ext cx=174514 cy=238200 (EMU). (The underlying custGeom path-point data
# for this freeform is not mutable through this host's COM surface, so only
# the shape's bounding box -- off/ext -- is updated here.)
$pl8 = $g.GroupItems.Item("pl8")
$pl8.Left   = 411.47984351968506
$pl8.Top    = 217.828031496063
$pl8.Width  = 13.741260842519685
$pl8.Height = 18.755906511811023

# "tx9" - "Plastid Parasite " label: off x=5014704 y=5359232 (EMU)
$tx9 = $g.GroupItems.Item("tx9")
$tx9.Left = 394.8585976771654
$tx9.Top  = 421.9867716535433

# "tx10" - " 83.721 %" label: off x=5309801 y=5665616 (EMU)
$tx10 = $g.GroupItems.Item("tx10")
$tx10.Left = 418.09456692913386
$tx10.Top  = 446.1114960629921

# "tx11" - "Non-Plastid Parasite " label: off x=3135160 y=3610502 (EMU)
$tx11 = $g.GroupItems.Item("tx11")
$tx11.Left = 246.86299912598426
$tx11.Top  = 284.2914960629921

# "tx12" - " 0.775 %" label: off x=3725512 y=3916886 (EMU)
$tx12 = $g.GroupItems.Item("tx12")
$tx12.Left = 293.34740157480314
$tx12.Top  = 308.41622047244095

# "tx13" - "Heterotroph " label: off x=4294554 y=2791736 (EMU)
$tx13 = $g.GroupItems.Item("tx13")
$tx13.Left = 338.15385826771654
$tx13.Top  = 219.82173928346458

# "tx14" - " 8.527 %" label: off x=4469230 y=3138666 (EMU)
$tx14 = $g.GroupItems.Item("tx14")
$tx14.Left = 351.907883015748
$tx14.Top  = 247.13906111811025

# "tx15" - "Endosymbiotic " label: off x=4751032 y=2230525 (EMU)
$tx15 = $g.GroupItems.Item("tx15")
$tx15.Left = 374.09700787401573
$tx15.Top  = 175.63188976377953

# "tx16" - " 6.977 %" label: off x=5058197 y=2579996 (EMU)
$tx16 = $g.GroupItems.Item("tx16")
$tx16.Left = 398.28322834645667
$tx16.Top  = 203.14929133858269
